{"js": "// Replace the date line and each '<a>\u00d7<b>=<c>' equation cell in the table\n// with its updated value, as a set of unique exact-text replacements.\nconst replacements = [\n  [\"2024-01-12 Friday\", \"2024-01-13 Saturday\"],\n  [\"626\u00d79=5634\", \"550\u00d79=4950\"],\n  [\"485\u00d76=2910\", \"612\u00d77=4284\"],\n  [\"304\u00d75=1520\", \"194\u00d75=970\"],\n  [\"656\u00d78=5248\", \"232\u00d76=1392\"],\n  [\"266\u00d73=798\", \"465\u00d75=2325\"],\n  [\"782\u00d75=3910\", \"533\u00d77=3731\"],\n  [\"177\u00d76=1062\", \"662\u00d75=3310\"],\n  [\"376\u00d77=2632\", \"737\u00d77=5159\"],\n  [\"715\u00d72=1430\", \"107\u00d75=535\"],\n  [\"613\u00d77=4291\", \"359\u00d75=1795\"],\n  [\"750\u00d76=4500\", \"475\u00d77=3325\"],\n  [\"542\u00d73=1626\", \"554\u00d72=1108\"],\n  [\"126\u00d72=252\", \"265\u00d77=1855\"],\n  [\"537\u00d75=2685\", \"738\u00d74=2952\"],\n  [\"559\u00d76=3354\", \"636\u00d79=5724\"],\n  [\"832\u00d76=4992\", \"393\u00d76=2358\"],\n  [\"500\u00d74=2000\", \"700\u00d73=2100\"],\n  [\"147\u00d77=1029\", \"241\u00d76=1446\"],\n  [\"952\u00d78=7616\", \"124\u00d78=992\"],\n  [\"457\u00d73=1371\", \"321\u00d74=1284\"],\n  [\"904\u00d79=8136\", \"418\u00d72=836\"],\n  [\"114\u00d74=456\", \"934\u00d75=4670\"],\n  [\"747\u00d73=2241\", \"416\u00d73=1248\"],\n  [\"994\u00d75=4970\", \"365\u00d76=2190\"],\n  [\"172\u00d74=688\", \"829\u00d77=5803\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each three-digit-by-one-digit multiplication\n# equation in the practice table to its new value (see commit diff).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-12 Friday\", \"2024-01-13 Saturday\"),\n    @(\"626\u00d79=5634\", \"550\u00d79=4950\"),\n    @(\"485\u00d76=2910\", \"612\u00d77=4284\"),\n    @(\"304\u00d75=1520\", \"194\u00d75=970\"),\n    @(\"656\u00d78=5248\", \"232\u00d76=1392\"),\n    @(\"266\u00d73=798\", \"465\u00d75=2325\"),\n    @(\"782\u00d75=3910\", \"533\u00d77=3731\"),\n    @(\"177\u00d76=1062\", \"662\u00d75=3310\"),\n    @(\"376\u00d77=2632\", \"737\u00d77=5159\"),\n    @(\"715\u00d72=1430\", \"107\u00d75=535\"),\n    @(\"613\u00d77=4291\", \"359\u00d75=1795\"),\n    @(\"750\u00d76=4500\", \"475\u00d77=3325\"),\n    @(\"542\u00d73=1626\", \"554\u00d72=1108\"),\n    @(\"126\u00d72=252\", \"265\u00d77=1855\"),\n    @(\"537\u00d75=2685\", \"738\u00d74=2952\"),\n    @(\"559\u00d76=3354\", \"636\u00d79=5724\"),\n    @(\"832\u00d76=4992\", \"393\u00d76=2358\"),\n    @(\"500\u00d74=2000\", \"700\u00d73=2100\"),\n    @(\"147\u00d77=1029\", \"241\u00d76=1446\"),\n    @(\"952\u00d78=7616\", \"124\u00d78=992\"),\n    @(\"457\u00d73=1371\", \"321\u00d74=1284\"),\n    @(\"904\u00d79=8136\", \"418\u00d72=836\"),\n    @(\"114\u00d74=456\", \"934\u00d75=4670\"),\n    @(\"747\u00d73=2241\", \"416\u00d73=1248\"),\n    @(\"994\u00d75=4970\", \"365\u00d76=2190\"),\n    @(\"172\u00d74=688\", \"829\u00d77=5803\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute(\n        $oldText,    # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap: wdFindContinue\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace: wdReplaceAll\n    )\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
